# Repull data, push all data, mean calculation
# Update column F (dSF) values for the rows where the repulled data differs
# from the originally stored value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -1
    9  = 4
    16 = 0
    20 = 0
    21 = 3
    24 = 2
    29 = 0
    44 = -2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
